$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.946.10"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.223.50"
$ws.Range("E3").Value = "  -0.85%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.98"
$ws.Range("E5").Value = "  +3.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.631"
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.21"
$ws.Range("E7").Value = "  +3.02%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.599"
$ws.Range("E9").Value = "  +8.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.71"
$ws.Range("E10").Value = "  +11.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0968"
$ws.Range("E11").Value = "  -2.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.32"
$ws.Range("E12").Value = "  -1.25%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.27"
$ws.Range("E13").Value = "  +7.60%  "

$ws.Range("E14").Value = "  -0.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.554.43"
$ws.Range("E15").Value = "  -0.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.02"
$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.880"
$ws.Range("E17").Value = "  +1.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.226.99"
$ws.Range("E18").Value = "  -0.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.845.98"
$ws.Range("E19").Value = "  -0.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0964"
$ws.Range("E20").Value = "  -0.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.24"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.78"
$ws.Range("E22").Value = "  -0.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.73"
$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  +2.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.08"
$ws.Range("E25").Value = "  +11.70%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.53"
$ws.Range("E26").Value = "  +15.46%  "

$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.17%  "

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.54"
$ws.Range("E28").Value = "  +3.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  -1.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.61"
$ws.Range("E30").Value = "  -0.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.86"
$ws.Range("E31").Value = "  +1.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.123"
$ws.Range("E32").Value = "  +0.80%  "

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.124"
$ws.Range("E33").Value = "  -1.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.55"
$ws.Range("E34").Value = "  +4.78%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0739"
$ws.Range("E35").Value = "  +2.69%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.71"
$ws.Range("E36").Value = "  +0.61%  "

$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.25"
$ws.Range("E37").Value = "  +15.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.03"
$ws.Range("E38").Value = "  +5.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0307"
$ws.Range("E39").Value = "  +9.29%  "

$ws.Range("E40").Value = "  -0.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.74"
$ws.Range("E41").Value = "  +26.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.94"
$ws.Range("E42").Value = "  +0.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.11"
$ws.Range("E43").Value = "  -3.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.207"
$ws.Range("E44").Value = "  +8.93%  "

$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.82"
$ws.Range("E45").Value = "  -2.66%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.75"
$ws.Range("E46").Value = "  -6.18%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.103"
$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.65"
$ws.Range("E48").Value = "  +1.09%  "

$ws.Range("E49").Value = "  +0.22%  "

$ws.Range("E50").Value = "  +5.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.38"
$ws.Range("E51").Value = "  +1.65%  "
